$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold, bordered, centered) from H1 onto
# the two new header cells before putting values in them.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns I (I0) and J (IF) for rows 2-8
$data = @{
    2 = @(1, 5)
    3 = @(1, 5)
    4 = @(1, 5)
    5 = @(1, 3)
    6 = @(6, 7)
    7 = @(4, 5)
    8 = @(9, 9)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
